# B6-PowerPoint.pptx edit
# 1) Re-point the three data tables (slides 14, 15, 16) from the custom
#    "Table_0" style to the built-in "Medium Style 2" table style.
# 2) The deck's applied theme ("Integral" / Red Violet) is swapped for the
#    default "Office Theme" colour palette. The two themes already share an
#    identical font scheme and format scheme, so only the 12 theme colours
#    need to change; they are exposed (and genuinely persisted) through
#    Slide.ThemeColorScheme, which all slides/layouts share via the single
#    slide master's theme part.

$p = $ppt.ActivePresentation

# --- 1) Table styles -------------------------------------------------------
$newTableStyle = "{B6BAD057-B781-4135-8840-4482D48C3D12}"

foreach ($slideIndex in 14, 15, 16) {
    $slide = $p.Slides.Item($slideIndex)
    foreach ($shp in $slide.Shapes) {
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle($newTableStyle)
        }
    }
}

# --- 2) Theme colours (Integral/Red Violet -> Office) ----------------------
$tcs = $p.Slides.Item(1).ThemeColorScheme

$tcs.Item(1).RGB  = 0          # dk1      000000
$tcs.Item(2).RGB  = 16777215   # lt1      FFFFFF
$tcs.Item(3).RGB  = 6968388    # dk2      44546A
$tcs.Item(4).RGB  = 15132391   # lt2      E7E6E6
$tcs.Item(5).RGB  = 13998939   # accent1  5B9BD5
$tcs.Item(6).RGB  = 3243501    # accent2  ED7D31
$tcs.Item(7).RGB  = 10855845   # accent3  A5A5A5
$tcs.Item(8).RGB  = 49407      # accent4  FFC000
$tcs.Item(9).RGB  = 12874308   # accent5  4472C4
$tcs.Item(10).RGB = 4697456    # accent6  70AD47
$tcs.Item(11).RGB = 12673797   # hlink    0563C1
$tcs.Item(12).RGB = 7491477    # folHlink 954F72
